# Outline_of_Python_Tutorials.docx - "adding to tutorial outline"
#
# The document currently ends with:
#   ... "Default Variables"
#   (blank paragraph)
#   (blank paragraph)
#   (a paragraph holding nothing but a manual line break)
#   (a final, empty "List Paragraph"-styled paragraph)
#
# We place the cursor right in front of that manual line break and type a
# manual page break, the "Python Tutorial #2" heading, and the new outline
# (Dictionaries / Tuples / Pickle). The pre-existing line break ends up
# trailing the new "III. Pickle" line, and the final paragraph is cleared
# down to a bare, unstyled empty paragraph.

$d = $word.ActiveDocument

# The paragraph that holds only the manual line break is the second-to-last
# paragraph in the document (the last paragraph is the trailing empty
# "List Paragraph" placeholder at the very end of the file).
$paraCount = $d.Paragraphs.Count
$breakParaIndex = $paraCount - 1
$breakPara = $d.Paragraphs.Item($breakParaIndex)

$insertPoint = $breakPara.Range
$insertPoint.Collapse(1)

$lines = @(
  "Python Tutorial #2",
  "I. Dictionaries",
  "",
  "    A. Dictionary Basics",
  "       i. creating dictionaries",
  "       ii. Accessing values",
  "       iii. changing values",
  "       iv. Adding items to a dictionary",
  "       v. removing items",
  "       vi. creating lists of keys and values",
  "       vii. iteration",
  "       viii. membership with 'in' and 'not in'",
  "       ix. copying dictionaries",
  "       x. merging dictionaries",
  "        ",
  "    B. Dictionary with lists as Values",
  "",
  "    C. Nested Dictionaries",
  "     ",
  "II. Tuples",
  "",
  "III. Pickle"
)

# Chr(12) = manual page break (Ctrl+Enter); "`r" = paragraph mark (Enter).
$blob = [char]12 + "`r" + ($lines -join "`r")
$insertPoint.Text = $blob

# The outline paragraphs ("I. Dictionaries" .. "III. Pickle") pick up the
# ListParagraph style with a 1080-twip (54pt) left indent, matching the
# rest of the outline's list items. The page-break paragraph and the new
# "Python Tutorial #2" heading paragraph keep the plain/default formatting
# they inherited from typing.
$outlineStartIndex = $breakParaIndex + 2
$outlineEndIndex = $outlineStartIndex + $lines.Count - 2
for ($i = $outlineStartIndex; $i -le $outlineEndIndex; $i++) {
    $p = $d.Paragraphs.Item($i)
    $p.Range.Style = "List Paragraph"
    $p.Format.LeftIndent = 54
}

# The final paragraph of the document (previously a ListParagraph/ind=2160
# placeholder) becomes a bare, unstyled empty paragraph.
$lastPara = $d.Paragraphs.Item($d.Paragraphs.Count)
$lastPara.Range.Style = "Normal"
$lastPara.Format.LeftIndent = 0

Write-Host "done"
